$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.152.09'
$ws.Range("E2").Value = '  +0.85%  '

$ws.Range("D3").Value = '1.750.98'
$ws.Range("E3").Value = '  +0.22%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").Value = "'237.55"
$ws.Range("E5").Value = '  +2.97%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.48%  '

$ws.Range("D7").Value = "'0.5279"
$ws.Range("E7").Value = '  +2.25%  '

$ws.Range("D8").Value = "'0.2814"
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = "'0.06195"
$ws.Range("E9").Value = '  +1.24%  '

$ws.Range("D10").Value = '1.747.90'
$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("E11").Value = '  +3.21%  '

$ws.Range("D12").Value = "'15.52"
$ws.Range("E12").Value = '  +0.42%  '

$ws.Range("D13").Value = "'0.6473"
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").Value = "'4.637"
$ws.Range("E14").Value = '  +2.56%  '

$ws.Range("D15").Value = "'78.67"
$ws.Range("E15").Value = '  +2.40%  '

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = '  +0.49%  '

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = '  +0.32%  '

$ws.Range("D18").Value = '26.056.41'
$ws.Range("E18").Value = '  +0.44%  '

$ws.Range("D19").Value = "'11.81"
$ws.Range("E19").Value = '  +2.66%  '

$ws.Range("D20").Value = "'0.000006751"
$ws.Range("E20").Value = '  +1.87%  '

$ws.Range("D21").Value = '1.970.83'
$ws.Range("E21").Value = '  +0.64%  '

$ws.Range("D22").Value = "'4.345"
$ws.Range("E22").Value = '  +5.35%  '

$ws.Range("D23").Value = "'8.779"

$ws.Range("D24").Value = "'5.246"
$ws.Range("E24").Value = '  +1.93%  '

$ws.Range("D25").Value = "'139.35"
$ws.Range("E25").Value = '  -0.31%  '

$ws.Range("D26").Value = "'1.520"
$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").Value = "'15.36"
$ws.Range("E27").Value = '  +1.83%  '

$ws.Range("D28").Value = "'1.817"
$ws.Range("E28").Value = '  -0.12%  '

$ws.Range("D29").Value = "'105.55"
$ws.Range("E29").Value = '  +2.40%  '

$ws.Range("D30").Value = "'0.08318"
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").Value = "'3.816"
$ws.Range("E31").Value = '  +5.11%  '

$ws.Range("D32").Value = "'3.658"
$ws.Range("E32").Value = '  +6.54%  '

$ws.Range("D33").Value = "'0.04614"
$ws.Range("E33").Value = '  +4.60%  '

$ws.Range("D34").Value = "'2.651"
$ws.Range("E34").Value = '  +0.86%  '

$ws.Range("D35").Value = "'1.019"
$ws.Range("E35").Value = '  +3.89%  '

$ws.Range("D36").Value = "'0.6353"
$ws.Range("E36").Value = '  +4.46%  '

$ws.Range("D37").Value = "'2.706"
$ws.Range("E37").Value = '  +0.65%  '

$ws.Range("D38").Value = "'0.01623"
$ws.Range("E38").Value = '  +3.49%  '

$ws.Range("D39").Value = "'1.980"
$ws.Range("E39").Value = '  +2.75%  '

$ws.Range("D40").Value = "'1.000"
$ws.Range("E40").Value = '  +0.58%  '

$ws.Range("E41").Value = '  +2.06%  '

$ws.Range("D42").Value = "'0.3943"
$ws.Range("E42").Value = '  +2.21%  '

$ws.Range("D43").Value = "'0.7552"
$ws.Range("E43").Value = '  +3.77%  '

$ws.Range("D44").Value = "'5.067"
$ws.Range("E44").Value = '  +2.92%  '

$ws.Range("D45").Value = "'0.1155"
$ws.Range("E45").Value = '  +3.94%  '

$ws.Range("D46").Value = "'6.392"
$ws.Range("E46").Value = '  +0.48%  '

$ws.Range("D47").Value = "'0.05355"
$ws.Range("E47").Value = '  -1.80%  '

$ws.Range("D48").Value = "'31.17"
$ws.Range("E48").Value = '  +4.64%  '

$ws.Range("D49").Value = "'54.65"
$ws.Range("E49").Value = '  +3.87%  '

$ws.Range("D50").Value = "'0.3485"
$ws.Range("E50").Value = '  +2.20%  '

$ws.Range("D51").Value = "'7.594"
$ws.Range("E51").Value = '  +0.70%  '
